# Update "Orders" sheet (sheet1): change F71, and append rows 72-79.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

# Helper: write a text value into a cell. Values that look like plain
# numbers ("10", "5", "3", ...) would otherwise be auto-converted to a
# real number by Excel, so for those we force a Text number format
# first, so they are stored as text instead.
function Set-TextValue($sheet, $row, $col, $value) {
    $cell = $sheet.Cells.Item($row, $col)
    if ($value -match '^-?\d+(\.\d+)?$') {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $value
}

# Existing row 71: Number column changes from 1 to 10
Set-TextValue $ws 71 6 "10"

# New row 72
Set-TextValue $ws 72 3 "224_折射_Reflex_Rosa rugosa Thunb._10stems"
Set-TextValue $ws 72 6 "5"

# New row 73
Set-TextValue $ws 73 3 "411_紫罗兰白_violet white_undefined_1bunch"
Set-TextValue $ws 73 6 "15"

# New row 74
Set-TextValue $ws 74 3 "505_紫罗兰紫_violet purple_undefined_1bunch"
Set-TextValue $ws 74 6 "10"

# New row 75
Set-TextValue $ws 75 3 "412_紫罗兰粉_violet pink_undefined_1bunch"
Set-TextValue $ws 75 6 "10"

# New row 76
Set-TextValue $ws 76 1 "3"
Set-TextValue $ws 76 3 "753_蝴蝶洋牡丹黄_butterfly  Ranunculus_undefined_1bunch"
Set-TextValue $ws 76 6 "10"

# New row 77
Set-TextValue $ws 77 3 "752_蝴蝶洋牡丹橙_butterfly  Ranunculus_undefined_1bunch"
Set-TextValue $ws 77 6 "5"

# New row 78
Set-TextValue $ws 78 3 "480_蝴蝶洋牡丹红_butterfly  Ranunculus_undefined_1bunch"
Set-TextValue $ws 78 6 "10"

# New row 79
Set-TextValue $ws 79 3 "585_洋牡丹红_undefined_undefined_1bunch"
Set-TextValue $ws 79 6 "10"

# Update "Summary" sheet (sheet2): G2 gets extra digits appended.
$ws2 = $wb.Worksheets.Item("Summary")
Set-TextValue $ws2 2 7 "0151540401033532151014713101491410105510115111082615151515151041595010201555510101551055151051010101010151016651010101051510101051010"
